# Update cryptos list on Sun Aug 11 19:37:02 UTC 2024 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns, and re-seat three rows whose
# coin ranking order changed (rows 32-34 and 50-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.657.98"
$ws.Range("E2").Value = "  -2.12%  "
$ws.Range("D3").Value = "2.600.13"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "512.95"
$ws.Range("E5").Value = "  -2.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.07"
$ws.Range("E6").Value = "  -5.50%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.561"
$ws.Range("E8").Value = "  -5.29%  "
$ws.Range("D9").Value = "2.601.69"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.24"
$ws.Range("E10").Value = "  -7.18%  "
$ws.Range("E11").Value = "  -2.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.336"
$ws.Range("E12").Value = "  -3.64%  "
$ws.Range("E13").Value = "  -1.01%  "
$ws.Range("D14").Value = "3.060.70"
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("D15").Value = "59.671.69"
$ws.Range("E15").Value = "  -2.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.88"
$ws.Range("E16").Value = "  -3.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000135"
$ws.Range("E17").Value = "  -3.93%  "
$ws.Range("D18").Value = "2.609.87"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.57"
$ws.Range("E19").Value = "  -3.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "339.00"
$ws.Range("E20").Value = "  -4.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.23"
$ws.Range("E21").Value = "  -3.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.02"
$ws.Range("E22").Value = "  -3.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.63"
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.412"
$ws.Range("E25").Value = "  -3.39%  "
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").Value = "  -5.43%  "
$ws.Range("D28").Value = "0.0₃0789"
$ws.Range("E28").Value = "  -6.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.90"
$ws.Range("E29").Value = "  -6.36%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.56"
$ws.Range("E31").Value = "  -2.82%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "150.28"
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.88"
$ws.Range("E33").Value = "  -7.42%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.67"
$ws.Range("E34").Value = "  -3.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.87"
$ws.Range("E35").Value = "  -7.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.897"
$ws.Range("E36").Value = "  -5.08%  "
$ws.Range("E37").Value = "  -6.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.58"
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.837"
$ws.Range("E39").Value = "  -1.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.40"
$ws.Range("E40").Value = "  -6.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.55"
$ws.Range("E41").Value = "  -6.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "283.03"
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.619"
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0988"
$ws.Range("E45").Value = "  -2.52%  "
$ws.Range("E46").Value = "  -3.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.96"
$ws.Range("E47").Value = "  -3.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.37"
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("E49").Value = "  -3.81%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.930.91"
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.54"
$ws.Range("E51").Value = "  -6.96%  "
